# Ajout de code pour export
# Rename the sheet, populate the header/"Cumul" columns and apply the
# reporting look (bold header, thin borders, yellow highlight on the
# figures) to the ENCAISSEMENT export sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet title -----------------------------------------------------
$ws.Name = "ENCAISSEMENT"

# --- Header row values -------------------------------------------------
$ws.Range("B1").Value = "25/04/2024"
$ws.Range("C1").Value = "Cumul Mois"
$ws.Range("D1").Value = "Cumul Année"

# --- Fill in the "Cumul Mois" column (was all zeros) --------------------
$ws.Range("C2").Value = 580048.29
$ws.Range("C3").Value = 21742.39
$ws.Range("C4").Value = 5043137.91
$ws.Range("C5").Value = 151987.06
$ws.Range("C6").Value = 242367.05

# --- Formatting ---------------------------------------------------------
# Bold header row
$ws.Range("A1:D1").Font.Bold = $true

# Thin border around every data cell (label + figures)
$ws.Range("A2:D6").Borders.LineStyle = 1

# Yellow highlight on the numeric figures only
$ws.Range("B2:D6").Interior.Color = 65535
